$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 139, shifting existing rows 139-199 down to 140-200
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new record's values
$ws.Cells.Item(139, 1).Value = 3
$ws.Cells.Item(139, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(139, 3).Value = "Coquimbo"
$ws.Cells.Item(139, 4).Value = 44917
$ws.Cells.Item(139, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(139, 5).Value = 5
$ws.Cells.Item(139, 6).Value = 100112052
$ws.Cells.Item(139, 7).Value = "Albahaca"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 110
$ws.Cells.Item(139, 11).Value = 5500
$ws.Cells.Item(139, 12).Value = 6000
$ws.Cells.Item(139, 13).Value = 5773
$ws.Cells.Item(139, 14).Value = "`$/docena de matas"
$ws.Cells.Item(139, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(139, 16).Value = 962
$ws.Cells.Item(139, 17).Value = 6
$ws.Cells.Item(139, 18).Value = "Hortaliza"
